$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v2")

# New header row (row 25): labels "Serial 64:" and "Avg" for each of the three blocks
$ws.Range("A25").Value = "Serial 64:"
$ws.Range("B25").Value = "Avg"
$ws.Range("F25").Value = "Serial 64:"
$ws.Range("G25").Value = "Avg"
$ws.Range("K25").Value = "Serial 64:"
$ws.Range("L25").Value = "Avg"

# Block 1 (columns A/B) - raw serial measurements + average formula
$ws.Range("A26").Value = 0.000518
$ws.Range("A27").Value = 0.000517
$ws.Range("A28").Value = 0.000557
$ws.Range("B26").Formula = "=(A26+A27+A28)/3"

# Block 2 (columns F/G)
$ws.Range("F26").Value = 0.024168
$ws.Range("F27").Value = 0.02402
$ws.Range("F28").Value = 0.024092
$ws.Range("G26").Formula = "=(F26+F27+F28)/3"

# Block 3 (columns K/L)
$ws.Range("K26").Value = 1.420498
$ws.Range("K27").Value = 1.421127
$ws.Range("K28").Value = 1.4210020000000001
$ws.Range("L26").Formula = "=(K26+K27+K28)/3"

# Update selection to reflect where the user ended up after editing
[void]$ws.Range("K29").Select()
